$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column C values to True/False booleans (rows 2-5), header label for C1
$ws.Range("C1").Value = "True / False"
$ws.Range("C2").Value = $true
$ws.Range("C3").Value = $false
$ws.Range("C4").Value = $false
$ws.Range("C5").Value = $true

# Update the selection to C4
$ws.Range("C4").Select() | Out-Null
